# "suppression de la date" - remove the "Dakar le 13-10-2015" paragraph
# (whole paragraph, including its paragraph mark) from the document.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Dakar le 13-10-2015") {
        $p.Range.Delete()
        break
    }
}
